$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Snippets" table currently spans A1:D165. Add two new rows describing
# the new Shape/TextFrame textbox snippets, expanding it to A1:D167.
$tbl = $ws.ListObjects.Item(1)

# Row 166: ShapeCollection.addTextBox -> excel-shape-textboxes / createTextbox
$row166 = $tbl.ListRows.Add()
$row166.Range.Cells.Item(1, 1).Value = "ShapeCollection"
$row166.Range.Cells.Item(1, 2).Value = "addTextBox"
$row166.Range.Cells.Item(1, 3).Value = "excel-shape-textboxes"
$row166.Range.Cells.Item(1, 4).Value = "createTextbox"

# Row 167: TextFrame.deleteText -> excel-shape-textboxes / deleteText
# (set column B/deleteText before column A/TextFrame so the shared-string
# table records "deleteText" right after "createTextbox", matching the
# order new strings were appended to the workbook)
$row167 = $tbl.ListRows.Add()
$row167.Range.Cells.Item(1, 2).Value = "deleteText"
$row167.Range.Cells.Item(1, 1).Value = "TextFrame"
$row167.Range.Cells.Item(1, 3).Value = "excel-shape-textboxes"
$row167.Range.Cells.Item(1, 4).Value = "deleteText"

# Move the view/selection to the newly added row, matching the saved
# worksheet view (scrolled so row 119 is visible, active cell on B167).
$win = $excel.ActiveWindow
$win.ScrollRow = 119
$win.ScrollColumn = 1
[void]$ws.Range("B167").Select()
